$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Daily Orders": a brand-new order (Order ID 8) came in, so it
# is inserted as the new top data row (row 2) and every existing
# order row shifts down by one.
# ------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Daily Orders")
$orders.Rows.Item(2).EntireRow.Insert()

$orders.Range("A2").Value = 8
$orders.Range("B2").Value = "2026-01-13 18:59"
$orders.Range("C2").Value = "Sagar Borse"
$orders.Range("D2").Value = "A-1608"
$orders.Range("E2").NumberFormat = "@"
$orders.Range("E2").Value = "7588930329"
$orders.Range("F2").Value = "Jawar Bhakari x1"
$orders.Range("G2").Value = 20
$orders.Range("H2").Value = "NEW"
$orders.Range("I2").Value = "PENDING"
$orders.Range("J2").NumberFormat = "@"
$orders.Range("J2").Value = "2026-01-16"
$orders.Range("K2").NumberFormat = "@"
$orders.Range("K2").Value = "10:00"
$orders.Range("L2:N2").NumberFormat = "@"
$orders.Range("L2").Value = ""
$orders.Range("M2").Value = ""
$orders.Range("N2").Value = ""

# ------------------------------------------------------------------
# Sheet "Summary": totals reflect the new order (+1 order, +1 new,
# +20 revenue).
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A2").Value = 8
$summary.Range("B2").Value = 7
$summary.Range("G2").Value = 215

# ------------------------------------------------------------------
# Sheet "Items Breakdown": "Jawar Bhakari" is a new item, inserted as
# row 4, pushing the existing "Onion Pakoda (Kanda Bhaje)" and "Pohe"
# rows down by one.
# ------------------------------------------------------------------
$items = $wb.Worksheets.Item("Items Breakdown")
$items.Rows.Item(4).EntireRow.Insert()
$items.Range("A4").Value = "Jawar Bhakari"
$items.Range("B4").Value = 1
$items.Range("C4").Value = 20
